$wb = $excel.ActiveWorkbook

# 1) Replace the "Ready for handoff" status text with "In Translation"
#    everywhere it occurs across all worksheets.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = [string]$cell.Value2
            if ("Ready for handoff" -eq $val) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# 2) Narrow the affected "status" columns to reflect the new, shorter text
#    (columns E & F on "Overview"; column C on "zh-cn" and "de-de").
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = 12.5

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = 12.5
